# Update cryptocurrency price and 1h-volume-change figures in the
# "cryptos" worksheet to reflect the latest scrape.
#
# Values are written with a leading apostrophe so that Excel stores them
# as text (matching the original inlineStr cells) instead of silently
# re-interpreting numeric-looking strings (e.g. "0.500", "16.30") as
# numbers and stripping significant trailing zeros, or mangling the
# multi-dot thousands-separated prices (e.g. "69.371.03").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.371.03"
$ws.Range("E2").Value = "'  -1.95%  "
$ws.Range("D3").Value = "'3.686.55"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'685.18"
$ws.Range("E5").Value = "'  -3.39%  "
$ws.Range("D6").Value = "'162.72"
$ws.Range("E6").Value = "'  -4.50%  "
$ws.Range("D7").Value = "'3.685.93"
$ws.Range("E7").Value = "'  -3.12%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "'  -4.12%  "
$ws.Range("E10").Value = "'  -7.38%  "
$ws.Range("D11").Value = "'7.26"
$ws.Range("E11").Value = "'  -2.45%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "'  -1.22%  "
$ws.Range("E13").Value = "'  -6.15%  "
$ws.Range("D14").Value = "'33.62"
$ws.Range("D15").Value = "'4.309.32"
$ws.Range("E15").Value = "'  -3.14%  "
$ws.Range("D16").Value = "'3.685.86"
$ws.Range("E16").Value = "'  -4.52%  "
$ws.Range("D17").Value = "'69.411.91"
$ws.Range("E17").Value = "'  -1.95%  "
$ws.Range("E18").Value = "'  -1.80%  "
$ws.Range("D19").Value = "'16.30"
$ws.Range("E19").Value = "'  -5.93%  "
$ws.Range("E20").Value = "'  -6.71%  "
$ws.Range("D21").Value = "'481.76"
$ws.Range("E21").Value = "'  -2.47%  "
$ws.Range("D22").Value = "'9.82"
$ws.Range("E22").Value = "'  -7.66%  "
$ws.Range("D23").Value = "'0.668"
$ws.Range("E23").Value = "'  -8.22%  "
$ws.Range("D24").Value = "'80.08"
$ws.Range("E24").Value = "'  -5.12%  "
$ws.Range("D25").Value = "'3.831.86"
$ws.Range("E25").Value = "'  -3.19%  "
$ws.Range("D26").Value = "'0.0000129"
$ws.Range("E26").Value = "'  -10.63%  "
$ws.Range("E27").Value = "'  +0.10%  "
$ws.Range("E28").Value = "'  -4.91%  "
$ws.Range("E29").Value = "'  -7.89%  "
$ws.Range("D30").Value = "'1.84"
$ws.Range("E30").Value = "'  -10.43%  "
$ws.Range("E31").Value = "'  -10.35%  "
$ws.Range("E32").Value = "'  -5.77%  "
$ws.Range("D33").Value = "'6.85"
$ws.Range("E33").Value = "'  -6.63%  "
$ws.Range("E34").Value = "'  -6.80%  "
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("D36").Value = "'0.164"
$ws.Range("E36").Value = "'  -4.90%  "
$ws.Range("D37").Value = "'3.649.58"
$ws.Range("E37").Value = "'  -3.34%  "
$ws.Range("E38").Value = "'  -5.70%  "
$ws.Range("D39").Value = "'6.08"
$ws.Range("E39").Value = "'  +2.44%  "
$ws.Range("D40").Value = "'0.0948"
$ws.Range("E40").Value = "'  -6.71%  "
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E42").Value = "'  -5.65%  "
$ws.Range("E43").Value = "'  -0.05%  "
$ws.Range("E44").Value = "'  -7.54%  "
$ws.Range("D45").Value = "'157.38"
$ws.Range("E45").Value = "'  -4.63%  "
$ws.Range("D46").Value = "'48.12"
$ws.Range("E46").Value = "'  -1.30%  "
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "'  -12.68%  "
$ws.Range("E48").Value = "'  -12.58%  "
$ws.Range("E49").Value = "'  -3.04%  "
$ws.Range("D50").Value = "'389.64"
$ws.Range("E50").Value = "'  -8.15%  "
$ws.Range("E51").Value = "'  -5.84%  "
